$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25; existing rows 25-30 shift down to 26-31.
$ws.Rows.Item(25).Insert()

# Fill the new row 25 with data, copying unchanged columns from the row
# that is now at 26 (the old row 25) and setting the new/changed values.
$ws.Cells.Item(25, 1).Value = 6
$ws.Cells.Item(25, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(25, 3).Value = "Metropolitana"
$ws.Cells.Item(25, 4).Value = 44627
$ws.Cells.Item(25, 4).NumberFormat = $ws.Cells.Item(26, 4).NumberFormat
$ws.Cells.Item(25, 5).Value = 13
$ws.Cells.Item(25, 6).Value = "Fruta"
$ws.Cells.Item(25, 7).Value = 100102
$ws.Cells.Item(25, 8).Value = "Cítricos"
$ws.Cells.Item(25, 9).Value = 100102006
$ws.Cells.Item(25, 10).Value = "Pomelo"
$ws.Cells.Item(25, 11).Value = "Start Ruby"
$ws.Cells.Item(25, 12).Value = "Primera"
$ws.Cells.Item(25, 13).Value = 6
$ws.Cells.Item(25, 14).Value = 240000
$ws.Cells.Item(25, 15).Value = 240000
$ws.Cells.Item(25, 16).Value = 240000
$ws.Cells.Item(25, 17).Value = "`$/bins (350 kilos)"
$ws.Cells.Item(25, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(25, 19).Value = 686
$ws.Cells.Item(25, 20).Value = 350
